$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename DeAndre Swift -> D'Andre Swift (row 10)
$ws.Range("A10").Value = "D'Andre Swift"

# Delete entire rows for players removed from rookie rights list
# Row numbers based on original layout; delete from bottom up to keep indices valid
$ws.Rows.Item(87).Delete()  # Jeremy Ruckert
$ws.Rows.Item(50).Delete()  # Dyami Brown
$ws.Rows.Item(49).Delete()  # Ian Book
$ws.Rows.Item(20).Delete()  # Lynn Bowden
$ws.Rows.Item(19).Delete()  # KJ Hamler

# Update the view state to match where the user ended up scrolled/selected
$ws.Range("D76").Select()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
